$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170793771743774
$ws.Range("B1").Value = 2.437577247619629
$ws.Range("D1").Value = 2.363499402999878
$ws.Range("E1").Value = 1.239295125007629
